$wb = $excel.ActiveWorkbook

# "assay_type list" sheet: fix the capitalization duplicate ("Publication" vs
# "publication") by keeping a single lowercase "publication" entry.
$ws3 = $wb.Worksheets.Item("assay_type list")
$ws3.Range("A1").Value = "publication"
$ws3.Range("A2").EntireRow.Delete()

# "Export as TSV" sheet: the assay_type column validation list now only has
# one valid value, so update the validation range and error message.
$ws1 = $wb.Worksheets.Item("Export as TSV")
$dv = $ws1.Range("B2:B1048576").Validation
$dv.Formula1 = "='assay_type list'!`$A`$1:`$A`$1"
$dv.ErrorMessage = "Value must be one of: publication."
